# UPDATE: data 25 Maret 2020 21:06 WIB
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new data row for 2020-03-25 first so the new footnote string
# is registered (and indexed) before the reworded old one.
$ws.Range("A38").Value = 43915
$ws.Range("A38").NumberFormat = "yyyy\-mm\-dd"
$ws.Range("B38").Value = 3332
$ws.Range("C38").Value = 790
$ws.Range("D38").Value = 31
$ws.Range("E38").Value = 58
$ws.Range("F38").Value = 2625
$ws.Range("G38").Value = 0
$ws.Range("H38").Value = 701
$ws.Range("I38").Value = "(informasi diambil di situs, pada 2020-03-25 21:06 WIB)"

# The old footnote text (row 37, col I) changes wording slightly
# ("pada situs" -> "di situs") and is kept as-is for the existing row.
$ws.Range("I37").Value = "(informasi kasus perawatan diambil di situs, pada 2020-03-24 16:07 WIB)"

# Update the view: scroll so row 10 is at top, select G38
$ws.Application.ActiveWindow.ScrollRow = 10
$ws.Range("G38").Select()
